$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wholeMatch = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# Rename province labels (column A) from mixed-case Indonesian names to
# uppercase names, using whole-cell matching so that substrings (e.g.
# "Riau" inside "Kepulauan Riau") are not partially replaced.
$ws.Cells.Replace("Aceh", "ACEH", $wholeMatch)
$ws.Cells.Replace("Sumatera Utara", "SUMATERA UTARA", $wholeMatch)
$ws.Cells.Replace("Sumatera Barat", "SUMATERA BARAT", $wholeMatch)
$ws.Cells.Replace("Riau", "RIAU", $wholeMatch)
$ws.Cells.Replace("Jambi", "JAMBI", $wholeMatch)
$ws.Cells.Replace("Sumatera Selatan", "SUMATERA SELATAN", $wholeMatch)
$ws.Cells.Replace("Bengkulu", "BENGKULU", $wholeMatch)
$ws.Cells.Replace("Lampung", "LAMPUNG", $wholeMatch)
$ws.Cells.Replace("Kepulauan Bangka Belitung", "KEPULAUAN BANGKA BELITUNG", $wholeMatch)
$ws.Cells.Replace("Kepulauan Riau", "KEPULAUAN RIAU", $wholeMatch)
$ws.Cells.Replace("DKI Jakarta", "DKI JAKARTA", $wholeMatch)
$ws.Cells.Replace("Jawa Barat", "JAWA BARAT", $wholeMatch)
$ws.Cells.Replace("Jawa Tengah", "JAWA TENGAH", $wholeMatch)
$ws.Cells.Replace("DI Yogyakarta", "DAERAH ISTIMEWA YOGYAKARTA", $wholeMatch)
$ws.Cells.Replace("Jawa Timur", "JAWA TIMUR", $wholeMatch)
$ws.Cells.Replace("Banten", "BANTEN", $wholeMatch)
$ws.Cells.Replace("Bali", "BALI", $wholeMatch)
$ws.Cells.Replace("Nusa Tenggara Barat", "NUSA TENGGARA BARAT", $wholeMatch)
$ws.Cells.Replace("Nusa Tenggara Timur", "NUSA TENGGARA TIMUR", $wholeMatch)
$ws.Cells.Replace("Kalimantan Barat", "KALIMANTAN BARAT", $wholeMatch)
$ws.Cells.Replace("Kalimantan Tengah", "KALIMANTAN TENGAH", $wholeMatch)
$ws.Cells.Replace("Kalimantan Selatan", "KALIMANTAN SELATAN", $wholeMatch)
$ws.Cells.Replace("Kalimantan Timur", "KALIMANTAN TIMUR", $wholeMatch)
$ws.Cells.Replace("Kalimantan Utara", "KALIMANTAN UTARA", $wholeMatch)
$ws.Cells.Replace("Sulawesi Utara", "SULAWESI UTARA", $wholeMatch)
$ws.Cells.Replace("Sulawesi Tengah", "SULAWESI TENGAH", $wholeMatch)
$ws.Cells.Replace("Sulawesi Selatan", "SULAWESI SELATAN", $wholeMatch)
$ws.Cells.Replace("Sulawesi Tenggara", "SULAWESI TENGGARA", $wholeMatch)
$ws.Cells.Replace("Gorontalo", "GORONTALO", $wholeMatch)
$ws.Cells.Replace("Sulawesi Barat", "SULAWESI BARAT", $wholeMatch)
$ws.Cells.Replace("Maluku", "MALUKU", $wholeMatch)
$ws.Cells.Replace("Maluku Utara", "MALUKU UTARA", $wholeMatch)
$ws.Cells.Replace("Papua Barat", "PAPUA BARAT", $wholeMatch)
$ws.Cells.Replace("Papua", "PAPUA", $wholeMatch)

# The original file applied a one-off Times New Roman style and a taller
# row height to the province column; the edited file reverts those cells
# back to the sheet's normal style/row height.
$ws.Range("A2:A137").Style = "Normal"
$ws.Rows("2:137").AutoFit()

# Reflect the scrolled viewport / selection the author left the sheet in.
$ws.Range("F138").Select()
